# "start of electronic box design"
# The "Misc" sheet (BOM data) had three "WIP" markers (column D) removed
# for rows 10, 13 and 14 -- rows 11 and 12 keep their "WIP" marker.
# The description column (C) was also manually narrowed from an
# autofit/best-fit width down to a fixed custom width, and the last
# selected cell on that sheet moved from C20 to A17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Misc")

# Misc is already the active sheet/tab in the workbook, make sure it stays so.
$ws.Activate()

# Remove the "WIP" note from rows 10, 13 and 14 (rows 11 & 12 keep theirs).
$ws.Range("D10").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("D14").ClearContents()

# Column C was resized by hand (no longer "best fit"/autofit).
$ws.Columns("C").ColumnWidth = 58.736979166666664

# Update the remembered selection on the sheet.
$ws.Range("A17").Select()
